$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct two values in the latest existing row (38) -- revised totals
$ws.Range("C38").Value = 1345101
$ws.Range("F38").Value = 366667

# Append the new day's row (39) with the Tompkins-county-only update;
# most national/state/county feeds are still pending so they are 0.
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("O39").Value = 574
$ws.Range("P39").Value = 103
$ws.Range("Q39").Value = 1456
$ws.Range("R39").Value = 2133
$ws.Range("S39").Value = 0
$ws.Range("T39").Value = 0
$ws.Range("U39").Value = 0
$ws.Range("V39").Value = 3
$ws.Range("W39").Value = 2
$ws.Range("X39").Value = 0
$ws.Range("Y39").Value = 69

# Carry over the bold/bordered "index" formatting used throughout column A
# (copy format only, so the freshly-written value in A39 is untouched).
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$excel.CutCopyMode = $false
